$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) B9:M10 spill block: _xll.MakeList() as an array formula over B9:M10 ---
# Seed the non-anchor spilled cells with their literal display values first
# (Excel stores only the top-left cell's formula for a CSE array, but every
# cell in the array range keeps its own cached value).
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 6
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 8
$ws.Range("K9").Value = 9
$ws.Range("L9").Value = 10
$ws.Range("M9").Value = 11

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 16
$ws.Range("G10").Value = 25
$ws.Range("H10").Value = 36
$ws.Range("I10").Value = 49
$ws.Range("J10").Value = 64
$ws.Range("K10").Value = 81
$ws.Range("L10").Value = 100
$ws.Range("M10").Value = 121

# Now enter the array formula over the whole block (anchor cell = B9)
$ws.Range("B9:M10").FormulaArray = "=_xll.MakeList()"

# --- 2) D11: a plain (non-array) formula ---
$ws.Range("D11").Formula = "=_xll.FooA()"

# --- 3) F14:Q15 spill block: _xll.MakeList() as an array formula over F14:Q15 ---
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 2
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 4
$ws.Range("K14").Value = 5
$ws.Range("L14").Value = 6
$ws.Range("M14").Value = 7
$ws.Range("N14").Value = 8
$ws.Range("O14").Value = 9
$ws.Range("P14").Value = 10
$ws.Range("Q14").Value = 11

$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = 25
$ws.Range("L15").Value = 36
$ws.Range("M15").Value = 49
$ws.Range("N15").Value = 64
$ws.Range("O15").Value = 81
$ws.Range("P15").Value = 100
$ws.Range("Q15").Value = 121

$ws.Range("F14:Q15").FormulaArray = "=_xll.MakeList()"

# --- 4) H20:S20 spill block: _xll.MakeArrayAndResize(1,12) as an array formula ---
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 4
$ws.Range("M20").Value = 5
$ws.Range("N20").Value = 6
$ws.Range("O20").Value = 7
$ws.Range("P20").Value = 8
$ws.Range("Q20").Value = 9
$ws.Range("R20").Value = 10
$ws.Range("S20").Value = 11

$ws.Range("H20:S20").FormulaArray = "=_xll.MakeArrayAndResize(1,12)"

# --- 5) Selection moves to H20 ---
$ws.Range("H20").Select()

# --- 6) Window position (best effort) ---
$excel.ActiveWindow.Left = 240
$excel.ActiveWindow.Top = 180
